# edit.ps1
# Applies the "act tablas web jul25" update to 420101.xlsx:
#  - Data sheet: adds a second data series ("DINEM - MIDES") in column C,
#    renames the old "Valor" column header to "MIDES-MEF-OPP", refreshes the
#    MIDES-MEF-OPP series values (2005-2023) and extends the years covered
#    down to 1985 / up to 2023.
#  - Metadata sheet: rewrites the "observaciones" text, and adds a new
#    "actualizacion" = "Julio 2025" row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) "Data" sheet
# ---------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")

$wsData.Cells.Item(1,1).Value = "Fecha"
$wsData.Cells.Item(1,2).Value = "MIDES-MEF-OPP"
$wsData.Cells.Item(1,3).Value = "DINEM - MIDES"

# Years go in column A as text (matching the original file's representation)
$wsData.Range("A2:A40").NumberFormat = "@"

$dataRows = @(
    @(2023, 6.5, $null),
    @(2022, 6, $null),
    @(2021, 6.8, $null),
    @(2020, 6.7, $null),
    @(2019, 6.3, $null),
    @(2018, 6.2, 6.7),
    @(2017, 6.2, 6.7),
    @(2016, 6, 6.5),
    @(2015, 5.8, 6.3),
    @(2014, 5.6, 6.1),
    @(2013, 5.5, 5.9),
    @(2012, 5.3, 5.7),
    @(2011, 4.9, 5.3),
    @(2010, 4.8, 5.2),
    @(2009, 4.7, 5.1),
    @(2008, 4.5, 4.9),
    @(2007, 3.8, 4.2),
    @(2006, 4, 4.3),
    @(2005, 3.7, 4),
    @(2004, $null, 3.2),
    @(2003, $null, 3.3),
    @(2002, $null, 3.5),
    @(2001, $null, 3.6),
    @(2000, $null, 3.4),
    @(1999, $null, 3.3),
    @(1998, $null, 3.1),
    @(1997, $null, 3.1),
    @(1996, $null, 3.2),
    @(1995, $null, 3.1),
    @(1994, $null, 3.1),
    @(1993, $null, 2.8),
    @(1992, $null, 2.7),
    @(1991, $null, 2.8),
    @(1990, $null, 2.7),
    @(1989, $null, 2.8),
    @(1988, $null, 2.6),
    @(1987, $null, 2.5),
    @(1986, $null, 2.7),
    @(1985, $null, 2.4)
)

$r = 2
foreach ($row in $dataRows) {
    $wsData.Cells.Item($r, 1).Value = [string]$row[0]

    if ($row[1] -ne $null) {
        $wsData.Cells.Item($r, 2).Value = $row[1]
    } else {
        $wsData.Cells.Item($r, 2).ClearContents()
    }

    if ($row[2] -ne $null) {
        $wsData.Cells.Item($r, 3).Value = $row[2]
    } else {
        $wsData.Cells.Item($r, 3).ClearContents()
    }

    $r = $r + 1
}

# ---------------------------------------------------------------
# 2) "Metadata" sheet
# ---------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

$wsMeta.Cells.Item(1,1).Value = " "

$wsMeta.Cells.Item(8,2).Value = "Las dos líneas representan metodologías ligeramente diferentes de cálculo. De acuerdo a lo informado en el Observatorio Social de MIDES, a partir del año 2016 se introdujo cambios en la metodología de estimación del Gasto Público Social producto de los cambios en la información brindada por el Presupuesto Nacional, lo cual llevó a trabajar en base al presupuesto por áreas programáticas (AP) de los incisos gubernamentales. El Gasto Público Social en Cultura y Deporte era considerado anteriormente bajo la denominación de Gasto Público Social No Convencional, definido como un subcomponente heterogéneo del GPS. La función Cultura y Deporte agrupa los gastos en museos, bibliotecas, organizaciones de prensa, servicios de televisión, deportes, y que antes también incluía otros conceptos que aludían a un aspecto multidisciplinario de los programas sociales. Se hizo una revisión de forma de dar consistencia en los conceptos para la serie desde 2015. La estimación siempre refiere a montos en pesos corrientes monto obligado intervenido por balance a partir de la información proporcionada mayoritariamente por Contaduría General de la Nación (CGN) del Ministerio de Economía y Finanzas (MEF). Para los años 2020 y 2021 se incluyen las erogaciones del fondo COVID destinadas a atender la emergencia sanitaria."

$wsMeta.Rows.Item(9).Insert()
$wsMeta.Cells.Item(9,1).Value = "actualizacion"
$wsMeta.Cells.Item(9,2).Value = "Julio 2025"

Write-Output "Edit complete"
